$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values stay text (many look numeric, e.g. "1.000", "106.24")
# so set the NumberFormat to Text ("@") before assigning, matching the source
# data which stores these as plain inline strings.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.951.63'
$ws.Range('E2').Value = '  -3.14%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.859.29'
$ws.Range('E3').Value = '  -2.46%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '318.18'
$ws.Range('E5').Value = '  -2.04%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.07%  '

$ws.Range('E7').Value = '  -4.67%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3700'
$ws.Range('E8').Value = '  -2.93%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07508'
$ws.Range('E9').Value = '  -2.66%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9367'
$ws.Range('E10').Value = '  -4.39%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '21.24'
$ws.Range('E11').Value = '  -3.96%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.896.65'
$ws.Range('E12').Value = '  +0.21%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.728'
$ws.Range('E13').Value = '  -3.18%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.442'
$ws.Range('E14').Value = '  -4.08%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.06840'
$ws.Range('E15').Value = '  -3.12%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  +0.00%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '81.64'
$ws.Range('E17').Value = '  -2.52%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000009040'
$ws.Range('E18').Value = '  -4.39%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.000'

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.95'
$ws.Range('E20').Value = '  -3.94%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '27.930.15'
$ws.Range('E21').Value = '  -3.25%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.111'
$ws.Range('E22').Value = '  -3.83%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.05'
$ws.Range('E23').Value = '  +1.23%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.102.46'
$ws.Range('E24').Value = '  -0.49%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.007'
$ws.Range('E25').Value = '  -4.33%  '

$ws.Range('E26').Value = '  -2.54%  '

$ws.Range('E27').Value = '  -3.34%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.418'
$ws.Range('E28').Value = '  -4.33%  '

$ws.Range('E29').Value = '  -3.50%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.739'
$ws.Range('E30').Value = '  -7.04%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08992'
$ws.Range('E31').Value = '  -3.11%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.8112'
$ws.Range('E32').Value = '  -6.22%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.821'
$ws.Range('E33').Value = '  -5.25%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.175'
$ws.Range('E34').Value = '  -5.52%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.943'
$ws.Range('E35').Value = '  -2.86%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  -0.08%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.05495'
$ws.Range('E37').Value = '  -3.76%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.113'
$ws.Range('E38').Value = '  -3.78%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01974'
$ws.Range('E39').Value = '  -3.22%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.896'
$ws.Range('E40').Value = '  +0.90%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.5258'
$ws.Range('E41').Value = '  -4.17%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '7.048'
$ws.Range('E42').Value = '  -5.39%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1689'
$ws.Range('E43').Value = '  -3.62%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.799'
$ws.Range('E44').Value = '  -5.69%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.06784'
$ws.Range('E45').Value = '  -1.39%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4904'
$ws.Range('E46').Value = '  -5.16%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.60'
$ws.Range('E47').Value = '  -5.36%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '106.24'
$ws.Range('E48').Value = '  -3.83%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.681'
$ws.Range('E49').Value = '  -5.42%  '

$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.000'
$ws.Range('E50').Value = '  -0.19%  '

$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.905'
$ws.Range('E51').Value = '  -11.99%  '

Write-Output "edits applied"